# Update gh-pages to output generated at 456a3b4
# Applies: bump a handful of "想去人数" (want-to-go) counts on the
# "展览" and "全部类型" sheets, and append a new row (row 18 / row 22
# respectively) for a newly-scraped event.

function Update-SheetCounts($ws, $newRow, $updates) {
    # $updates: F-column count bumps, keyed by row number on this sheet
    # (the two sheets list the same events in different row order, so
    # each caller supplies its own row->newvalue map).
    foreach ($r in $updates.Keys) {
        $ws.Cells.Item($r, 6).Value = $updates[$r]
    }

    # Append the new event row. Column A uses the same bold/centered/
    # thin-bordered look as every other index cell in the sheet.
    $idxCell = $ws.Cells.Item($newRow, 1)
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108  # xlCenter
    $idxCell.VerticalAlignment = -4160    # xlTop
    $idxCell.Borders.LineStyle = 1        # xlContinuous (thin, all sides)
    $idxCell.Value = $newRow - 1
    # Force the date-column cell to Text first so Excel doesn't coerce
    # this "yyyy-mm-dd" look-alike into a date serial — the source file
    # stores it as a plain string, like every other row in column B.
    $dateCell = $ws.Cells.Item($newRow, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2024-09-07"
    $ws.Cells.Item($newRow, 3).Value = "合肥·国乙only宇宙心动（含夜场）"
    $ws.Cells.Item($newRow, 4).Value = "文忠路1865号 赫拉诺言艺术中心"
    $ws.Cells.Item($newRow, 5).Value = "2024.09.07 10:00-09.07 21:00"
    $ws.Cells.Item($newRow, 6).Value = 6
    $ws.Cells.Item($newRow, 7).Value = 48
    $ws.Cells.Item($newRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89803"
    $ws.Cells.Item($newRow, 9).Value = "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"
}

$wb = $excel.ActiveWorkbook

# "展览" sheet: rows 3,4,6,9,10,11,13,14,15,16,17.
$updatesExhibit = @{
    3  = 553
    4  = 199
    6  = 508
    9  = 46
    10 = 6768
    11 = 235
    13 = 3073
    14 = 200
    15 = 352
    16 = 259
    17 = 549
}
$wsExhibit = $wb.Worksheets.Item("展览")
Update-SheetCounts $wsExhibit 18 $updatesExhibit

# "全部类型" sheet: same events, different row numbers (5,6,8,11,13,
# 15,17,18,19,20,21).
$updatesAll = @{
    5  = 553
    6  = 199
    8  = 508
    11 = 46
    13 = 6768
    15 = 235
    17 = 3073
    18 = 200
    19 = 352
    20 = 259
    21 = 549
}
$wsAll = $wb.Worksheets.Item("全部类型")
Update-SheetCounts $wsAll 22 $updatesAll
